$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
# New sale registered for ALMEIDA CUATIN JHONATHANN CARLOS / ALTAMIRANO MARCATOMA EDISON PAULINO
# in the PORCELANATO category -> updates the value and the "X de 30" summary count.
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M2").Value = 3548.83
$wsGrupo.Range("M32").Value = "8 de 30"

# --- Sheet 2: VENTA MENSUAL ---
# Same new sale booked in julio (column F) for the same advisor/client, plus the
# recalculated column total.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F2").Value = 4360.63
$wsMensual.Range("F32").Value = 12431.9

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
# PORCELANATO row (15) and TOTAL row (18) reflect the extra sale: VENTA goes up,
# POR CUMPLIR goes down, CUMPLIMIENTO (%) goes up.
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D15").Value = 9405.549999999999
$wsCumplimiento.Range("E15").Value = 14053.27
$wsCumplimiento.Range("F15").Value = 0.4009387513949977

$wsCumplimiento.Range("D18").Value = 12421.62
$wsCumplimiento.Range("E18").Value = 21513.09607548726
$wsCumplimiento.Range("F18").Value = 0.366044612613475
